# Keyboard Shortcuts workbook update
# - Adds a new keyboard-shortcut row to the "Keyboard Shortcuts" sheet describing
#   the new OPTION+B [Mac] / CTRL+SHIFT+B [Win] + <value> shortcut that fires a
#   custom event button action.
# - Updates the sheet's selection/active cell bookkeeping to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Keyboard Shortcuts")
$ws.Activate()

# New row 34: the OPTION+B / CTRL+SHIFT+B shortcut and its description.
# Single-quoted strings so PowerShell does not interpret the literal `\n`
# (the workbook stores these multi-shortcut cells with a literal backslash-n,
# not an actual line break, matching the existing rows such as A7/A39/A41).
$ws.Cells.Item(34, 1).Value = 'OPTION+B + <value> [Mac]\nCTRL+SHIFT+B + <value> [Win]'
$ws.Cells.Item(34, 2).Value = 'Fire custom event button action. Value is a two digit number indicating the button number.'

# Match the author's row-height tweak on the row just above the new entry.
$ws.Rows.Item(33).RowHeight = 13.8

# Update the sheet's selection to the new active cell.
$ws.Range("B33").Select() | Out-Null
